# Schulferien - unnötige Zeilen gelöscht
# Clear the leftover "footnote" block (rows 20-26) that held the
# "Stand / Angegeben ist ... / SH - Auf den Inseln ..." remarks. The
# row/column layout, styling and merged cells stay untouched - only the
# cell contents are removed, just like Delete/ClearContents in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20:G26").ClearContents()

# Reflect the new selection left behind in the sheet (as seen after
# selecting the now-empty block in the Excel UI).
$ws.Range("A18:G28").Select()
